$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.678.93'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').Value = '2.667.78'
$ws.Range('E3').Value = '  +4.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.590'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').Value = '2.666.60'
$ws.Range('E9').Value = '  +4.05%  '
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.71'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.07%  '
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.357'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.24%  '
$ws.Range('D15').Value = '3.141.07'
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('D16').Value = '63.533.93'
$ws.Range('E16').Value = '  +1.63%  '
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('D18').Value = '2.654.35'
$ws.Range('E18').Value = '  +3.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '341.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.41'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.81'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.50%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.69'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.166'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.99%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '540.86'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +18.69%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.91'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  +15.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.61%  '
$ws.Range('D34').Value = '0.0₃0817'
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '174.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +17.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.405'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.83'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '173.30'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.78'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.12%  '
$ws.Range('E46').Value = '  +6.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.638'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.55%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0241'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.42%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0964'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.83'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.74'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.53%  '
